# Update LR-pairs table with new TPM-based values (adds "Resolving-Mac" cluster)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colNum = @{
    A = 1
    B = 2
    C = 3
    D = 4
    E = 5
    F = 6
    G = 7
    H = 8
    I = 9
    J = 10
    K = 11
    L = 12
    M = 13
    N = 14
    O = 15
    P = 16
    Q = 17
    R = 18
    S = 19
    T = 20
}

$rowData = @(
    @{ row=2; A="ECs"; B="Agt"; C="Agtr2"; D="ECs"; E=1; F=0.3333333333333333; G=0.05804633333333333; H=0.174139; I=0.02760497488446473; J=0.02760497488446473; K=1; L=0.3333333333333333; M=0.045282; N=0.135846; O=0.001264082841858775; P=0.001264082841858776; Q=0.002628454066; R=0.023656086594; S=0.0000348949751013943; T=0.0000348949751013943 },
    @{ row=3; A="ECs"; B="Agt"; C="Agtr2"; D="FAPs"; E=1; F=0.3333333333333333; G=0.05804633333333333; H=0.174139; I=0.02760497488446473; J=0.02760497488446473; K=3; L=1; M=34.682839; N=104.048517; O=0.9681988800594137; P=0.9681988800594138; Q=2.013211633540333; R=18.118904701863; S=0.026727105767207; T=0.026727105767207 },
    @{ row=4; A="ECs"; B="Agt"; C="Agtr2"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.05804633333333333; H=0.174139; I=0.02760497488446473; J=0.02760497488446473; K=2; L=0.6666666666666666; M=1.093898333333333; N=3.281695; O=0.03053703709872749; P=0.03053703709872749; Q=0.06349678728944445; R=0.5714710856049999; S=0.00084297414215634; T=0.0008429741421563399 },
    @{ row=5; A="FAPs"; B="Agt"; C="Agtr2"; D="ECs"; E=3; F=1; G=1.761113666666667; H=5.283341; I=0.8375291899635509; J=0.8375291899635509; K=1; L=0.3333333333333333; M=0.045282; N=0.135846; O=0.001264082841858775; P=0.001264082841858776; Q=0.07974674905399999; R=0.717720741486; S=0.001058706278588803; T=0.001058706278588804 },
    @{ row=6; A="FAPs"; B="Agt"; C="Agtr2"; D="FAPs"; E=3; F=1; G=1.761113666666667; H=5.283341; I=0.8375291899635509; J=0.8375291899635509; K=3; L=1; M=34.682839; N=104.048517; O=0.9681988800594137; P=0.9681988800594138; Q=61.08042176169967; R=549.723795855297; S=0.8108948237397778; T=0.8108948237397779 },
    @{ row=7; A="FAPs"; B="Agt"; C="Agtr2"; D="MuSCs"; E=3; F=1; G=1.761113666666667; H=5.283341; I=0.8375291899635509; J=0.8375291899635509; K=2; L=0.6666666666666666; M=1.093898333333333; N=3.281695; O=0.03053703709872749; P=0.03053703709872749; Q=1.926479304777222; R=17.338313742995; S=0.02557565994518413; T=0.02557565994518413 },
    @{ row=8; A="MuSCs"; B="Agt"; C="Agtr2"; D="ECs"; E=2; F=0.6666666666666666; G=0.2331273333333333; H=0.6993819999999999; I=0.1108678845327394; J=0.1108678845327394; K=1; L=0.3333333333333333; M=0.045282; N=0.135846; O=0.001264082841858775; P=0.001264082841858776; Q=0.010556471908; R=0.09500824717199999; S=0.0001401461905510158; T=0.0001401461905510158 },
    @{ row=9; A="MuSCs"; B="Agt"; C="Agtr2"; D="FAPs"; E=2; F=0.6666666666666666; G=0.2331273333333333; H=0.6993819999999999; I=0.1108678845327394; J=0.1108678845327394; K=3; L=1; M=34.682839; N=104.048517; O=0.9681988800594137; P=0.9681988800594138; Q=8.085517768499333; R=72.76965991649399; S=0.1073421616391547; T=0.1073421616391547 },
    @{ row=10; A="MuSCs"; B="Agt"; C="Agtr2"; D="MuSCs"; E=2; F=0.6666666666666666; G=0.2331273333333333; H=0.6993819999999999; I=0.1108678845327394; J=0.1108678845327394; K=2; L=0.6666666666666666; M=1.093898333333333; N=3.281695; O=0.03053703709872749; P=0.03053703709872749; Q=0.2550176013877778; R=2.29515841249; S=0.003385576703033699; T=0.003385576703033699 },
    @{ row=11; A="Resolving-Mac"; B="Agt"; C="Agtr2"; D="ECs"; E=1; F=0.3333333333333333; G=0.05046166666666666; H=0.151385; I=0.02399795061924493; J=0.02399795061924493; K=1; L=0.3333333333333333; M=0.045282; N=0.135846; O=0.001264082841858775; P=0.001264082841858776; Q=0.002285005189999999; R=0.02056504671; S=0.00003033539761756169; T=0.00003033539761756169 },
    @{ row=12; A="Resolving-Mac"; B="Agt"; C="Agtr2"; D="FAPs"; E=1; F=0.3333333333333333; G=0.05046166666666666; H=0.151385; I=0.02399795061924493; J=0.02399795061924493; K=3; L=1; M=34.682839; N=104.048517; O=0.9681988800594137; P=0.9681988800594138; Q=1.750153860671666; R=15.751384746045; S=0.02323478891327405; T=0.02323478891327406 },
    @{ row=13; A="Resolving-Mac"; B="Agt"; C="Agtr2"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.05046166666666666; H=0.151385; I=0.02399795061924493; J=0.02399795061924493; K=2; L=0.6666666666666666; M=1.093898333333333; N=3.281695; O=0.03053703709872749; P=0.03053703709872749; Q=0.05519993306388889; R=0.496799397575; S=0.0007328263083533127; T=0.0007328263083533127 }
)

foreach ($r in $rowData) {
    $ws.Cells.Item($r.row, $colNum.A).Value = $r.A
    $ws.Cells.Item($r.row, $colNum.B).Value = $r.B
    $ws.Cells.Item($r.row, $colNum.C).Value = $r.C
    $ws.Cells.Item($r.row, $colNum.D).Value = $r.D
    $ws.Cells.Item($r.row, $colNum.E).Value = $r.E
    $ws.Cells.Item($r.row, $colNum.F).Value = $r.F
    $ws.Cells.Item($r.row, $colNum.G).Value = $r.G
    $ws.Cells.Item($r.row, $colNum.H).Value = $r.H
    $ws.Cells.Item($r.row, $colNum.I).Value = $r.I
    $ws.Cells.Item($r.row, $colNum.J).Value = $r.J
    $ws.Cells.Item($r.row, $colNum.K).Value = $r.K
    $ws.Cells.Item($r.row, $colNum.L).Value = $r.L
    $ws.Cells.Item($r.row, $colNum.M).Value = $r.M
    $ws.Cells.Item($r.row, $colNum.N).Value = $r.N
    $ws.Cells.Item($r.row, $colNum.O).Value = $r.O
    $ws.Cells.Item($r.row, $colNum.P).Value = $r.P
    $ws.Cells.Item($r.row, $colNum.Q).Value = $r.Q
    $ws.Cells.Item($r.row, $colNum.R).Value = $r.R
    $ws.Cells.Item($r.row, $colNum.S).Value = $r.S
    $ws.Cells.Item($r.row, $colNum.T).Value = $r.T
}

Write-Output ("Dimension=" + $ws.UsedRange.Address())
